$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlLeft = -4131 ; xlPasteFormats = -4122
$xlLeft = -4131
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# New user rows (11-15). Username/Password columns get a left-aligned
# "Normal" style, the Email column gets a left-aligned "Hyperlink" style
# (mirroring the existing C2:C10 hyperlink cells), and CreatedDate/
# LastAccessed keep the existing short-date number format (copied from the
# row above so the style id matches exactly).
# ---------------------------------------------------------------------------

function Add-UserRow($Row, $Username, $Email) {
    $prevRow = $Row - 1

    $ws.Range("A" + $Row).Value = $Username
    $ws.Range("B" + $Row).Value = $Username
    $ws.Range("A" + $Row + ":B" + $Row).HorizontalAlignment = $xlLeft

    $ws.Range("D" + $prevRow + ":E" + $prevRow).Copy()
    $ws.Range("D" + $Row + ":E" + $Row).PasteSpecial($xlPasteFormats)
    $ws.Range("D" + $Row).Value = 44075
    $ws.Range("E" + $Row).Value = 44083

    $ws.Range("F" + $Row).Value = "Sowmya"

    $ws.Range("C" + $Row).Value = $Email
    $ws.Range("C" + $Row).HorizontalAlignment = $xlLeft
    $ws.Hyperlinks.Add($ws.Range("C" + $Row), "mailto:" + $Email)
}

# Row order mirrors the original authoring session: 13, 12, 11, 14, 15.
Add-UserRow 13 "Dali" "Dali@gmail.com"
Add-UserRow 12 "Daley" "Daley@gmail.com"
Add-UserRow 11 "Curie" "Curie@gmail.com"
Add-UserRow 14 "Gary" "Gary@gmail.com"
Add-UserRow 15 "Emilysmith" "emily@gmail.com"

# Trailing blank row 16 - only E16 carries the date-format style forward.
$ws.Range("E15").Copy()
$ws.Range("E16").PasteSpecial($xlPasteFormats)

# Column widths (closest values reachable through the ColumnWidth setter).
$ws.Columns.Item(1).ColumnWidth = 14.6
$ws.Columns.Item(2).ColumnWidth = 15.5
$ws.Columns.Item(3).ColumnWidth = 22.5

# Selection moves to B2 in the saved view.
$ws.Range("B2").Select()

$excel.CutCopyMode = $false
